# edit.ps1 -- apply the "Proyecto/Tecnico MATRIX" title-retranslation edit
# described by the commit:
#   "This is a new comment / Resolved Jira # 1234"
#
# Changes made to the document (see diff):
#   1. The title line "PROJECT/TECHNICAL MATRIX!" becomes four separately
#      spell-checked runs: "Proyecto" / "/" / "Tecnico" (flagged with a
#      proofErr spellStart/spellEnd pair, since it is not a real Spanish
#      word) / " MATRIX", followed by the existing "!" run.
#   2. The "_GoBack" bookmark (which Word re-drops at the position of the
#      user's last edit) moves from the end of the title paragraph to the
#      end of the blank paragraph that precedes it.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Locate the title text. We search instead of assuming a fixed paragraph
# index so the script is resilient to minor structural differences.
# ---------------------------------------------------------------------
$titleRange = $d.Content
$found = $titleRange.Find.Execute("PROJECT/TECHNICAL MATRIX!", $true, $false,
                                   $false, $false, $false, $true, 1, $false,
                                   "", 0)
if (-not $found) {
    throw "Could not find the 'PROJECT/TECHNICAL MATRIX!' title run to edit."
}
$titleStart = $titleRange.Start
$titleEnd = $titleRange.End

# ---------------------------------------------------------------------
# Step 1: move the "_GoBack" bookmark from the end of the title paragraph
# to the end of the paragraph immediately before it (an empty paragraph).
# Re-adding a bookmark with the same name moves it -- Word only allows one
# bookmark per name -- so this both removes the old one and creates the
# new one.
# ---------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$prevPara = $null
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if (($p.Range.Start -le $titleStart) -and ($titleStart -lt $p.Range.End)) {
        if ($i -gt 1) {
            $prevPara = $d.Paragraphs.Item($i - 1)
        }
        break
    }
}

if ($prevPara -ne $null) {
    $bookmarkRange = $d.Range($prevPara.Range.Start, $prevPara.Range.Start)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)
}

# ---------------------------------------------------------------------
# Step 2: rewrite the title run ("PROJECT/TECHNICAL MATRIX!") as:
#   "Proyecto" + "/" + (spellStart)"Tecnico"(spellEnd) + " MATRIX" + "!"
# all bold, "Microsoft New Tai Lue", sz 28 / szCs 18 -- matching the
# formatting of the original runs -- with a proofErr pair flagging
# "Tecnico" as a misspelling (no accent) the way Word's spell checker
# would.
# ---------------------------------------------------------------------
$titleTextRange = $d.Range($titleStart, $titleEnd)

$rPr = '<w:rPr><w:rFonts w:ascii="Microsoft New Tai Lue" w:hAnsi="Microsoft New Tai Lue" w:cs="Microsoft New Tai Lue"/><w:b/><w:sz w:val="28"/><w:szCs w:val="18"/></w:rPr>'

$newRunsXml = (
    '<w:r>' + $rPr + '<w:t>Proyecto</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>/</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPr + '<w:t>Tecnico</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> MATRIX</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>!</w:t></w:r>'
)

$pPr = '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Microsoft New Tai Lue" w:hAnsi="Microsoft New Tai Lue" w:cs="Microsoft New Tai Lue"/><w:sz w:val="24"/><w:szCs w:val="18"/></w:rPr></w:pPr>'

$packageXml = (
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $pPr + $newRunsXml + '</w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
)

$titleTextRange.InsertXML($packageXml)

Write-Output "Title paragraph rewritten and _GoBack bookmark relocated."
